$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New source codes in column T for rows 17-19
$ws.Range("T17").Value = "AB12345"
$ws.Range("T18").Value = "CD56789"
$ws.Range("T19").Value = "EF01234"

# Dynamic-array formula: split each T value into two chunks (chars 1-2 and
# chars 6-7) using MID + CHOOSECOLS, spilling across V17:W19
$ws.Range("V17:W19").FormulaArray = "=CHOOSECOLS(MID(T17:T19,{1,3,6},{2,3,2}),1,3)"

# Move the active selection to reflect where the author ended up working
$ws.Range("V18").Select() | Out-Null
